$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update first-table input values (F/G columns); H column formulas auto-recalculate.
$ws.Range("F3").Value = 8.6739999999999995
$ws.Range("G3").Value = 6.1929999999999801

$ws.Range("F4").Value = 17.846
$ws.Range("G4").Value = 3.1840000000000201

$ws.Range("F5").Value = 46.921999999999997
$ws.Range("G5").Value = 32.091000000000001

$ws.Range("F6").Value = 220.60599999999999
$ws.Range("G6").Value = 19.498999999999999

# Strip the bold/black-font style from the first table's header row and
# row-label column so these cells fall back to the default "Normal" style.
$ws.Range("A2:C2").Style = "Normal"
$ws.Range("A3").Style = "Normal"
$ws.Range("A4").Style = "Normal"
$ws.Range("A5").Style = "Normal"
$ws.Range("A6").Style = "Normal"

# Move the active selection/cursor as recorded in the saved view state.
$ws.Range("H17").Select()
